$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The BOM was regenerated with one fewer component (the J1 / AVR_ICSP_3x2 /
# AVR-ISP-6 ICSP header row). Clear that row's contents in place -- the rows
# below it keep their original row numbers (this was not a row delete/shift).
$ws.Range("A12:G12").ClearContents()

# Leave the selection where the author last left it when saving.
[void]$ws.Range("C17").Select()
